# Add "Denver CO, Minneapolis MN, " in the middle of the run that reads
# "  The new states containing Seattle WA, San Francisco CA, and the urban
# corridor", splitting that single run into three runs:
#   1. "  The new states containing Seattle WA, San Francisco CA, "
#   2. "Denver CO, Minneapolis MN, "
#   3. "and the urban corridor"
#
# A plain Find/Replace (or InsertBefore/InsertAfter) applied directly on
# a Range inside that paragraph causes this engine to re-coalesce every
# run of the touched paragraph into a single run, which would lose the
# 3-way run split the target OOXML requires. Splitting the paragraph in
# two instead (which *does* preserve the run boundary at the split
# point), inserting the new text at the start of the new paragraph, and
# then deleting the paragraph mark to re-join the two paragraphs leaves
# all the run boundaries intact.
#
# Note: positions captured before a structural edit (like
# InsertParagraphAfter) can go stale, so every step below re-locates its
# anchor text with a fresh Find immediately before using it.

$d = $word.ActiveDocument

# 1) Find "San Francisco CA, " - split the paragraph right after it so
#    the existing run boundary there is preserved.
$find1 = $d.Content
$found1 = $find1.Find.Execute(
    "San Francisco CA, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'San Francisco CA, ' in the document"
}
$splitPos = $find1.End
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# 2) The new paragraph now begins with "and the urban corridor...".
#    Re-find it fresh (the paragraph split moved things around) and
#    insert the new text right before it. InsertBefore on a collapsed
#    range creates its own run instead of merging into the following
#    run.
$find2 = $d.Content
$found2 = $find2.Find.Execute(
    "and the urban corridor", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'and the urban corridor' after the paragraph split"
}
$insertPos = $find2.Start
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore("Denver CO, Minneapolis MN, ")

# 3) Re-join the two paragraphs by deleting the paragraph mark that sits
#    right after "San Francisco CA, " (re-find it fresh once more for an
#    accurate position). The runs on either side of the old split point
#    stay separate, giving the desired 3-run split.
$find3 = $d.Content
$found3 = $find3.Find.Execute(
    "San Francisco CA, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not re-find 'San Francisco CA, ' before merging paragraphs"
}
$mergePos = $find3.End
$markRange = $d.Range($mergePos, $mergePos + 1)
$markRange.Delete()

Write-Output "Done"
